$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '68.282.51'
$ws.Cells.Item(2, 5).Value = '  +2.93%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.638.94'
$ws.Cells.Item(3, 5).Value = '  +2.17%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.35%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '197.38'
$ws.Cells.Item(5, 5).Value = '  +10.06%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '578.67'
$ws.Cells.Item(6, 5).Value = '  -0.76%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '3.632.61'
$ws.Cells.Item(7, 5).Value = '  +2.29%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +2.47%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  -0.34%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  +1.80%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.155'
$ws.Cells.Item(11, 5).Value = '  +8.57%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '56.77'
$ws.Cells.Item(12, 5).Value = '  +6.58%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  +17.30%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '10.10'
$ws.Cells.Item(14, 5).Value = '  +2.70%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '4.220.03'
$ws.Cells.Item(15, 5).Value = '  +1.84%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '3.642.57'
$ws.Cells.Item(16, 5).Value = '  +2.33%  '

# Row 17
$ws.Cells.Item(17, 5).Value = '  +0.80%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '12.59'
$ws.Cells.Item(18, 5).Value = '  +4.37%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '68.299.54'
$ws.Cells.Item(19, 5).Value = '  +3.38%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '18.63'
$ws.Cells.Item(20, 5).Value = '  +2.46%  '

# Row 21
$ws.Cells.Item(21, 5).Value = '  +4.03%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '403.85'
$ws.Cells.Item(22, 5).Value = '  +3.53%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '13.19'
$ws.Cells.Item(23, 5).Value = '  +30.11%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  +0.34%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +2.40%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +4.48%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '12.67'
$ws.Cells.Item(27, 5).Value = '  +4.51%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '3.88'
$ws.Cells.Item(28, 5).Value = '  +7.90%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +1.68%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '8.27'
$ws.Cells.Item(30, 5).Value = '  +23.56%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +3.97%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '31.84'
$ws.Cells.Item(32, 5).Value = '  +2.99%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '688.95'
$ws.Cells.Item(33, 5).Value = '  +16.83%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  +3.77%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  +5.59%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '64.66'
$ws.Cells.Item(36, 5).Value = '  -0.82%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '42.83'
$ws.Cells.Item(37, 5).Value = '  +4.33%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.429'
$ws.Cells.Item(38, 5).Value = '  +16.23%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  +0.10%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  +8.45%  '

# Row 41
$ws.Cells.Item(41, 2).Value = 'Kaspa'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.140'
$ws.Cells.Item(41, 5).Value = '  +8.96%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'Fetch.AI'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '2.90'
$ws.Cells.Item(42, 5).Value = '  +21.59%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '3.18'
$ws.Cells.Item(43, 5).Value = '  +15.76%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '3.222.01'
$ws.Cells.Item(44, 5).Value = '  +17.30%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '3.08'
$ws.Cells.Item(45, 5).Value = '  +42.72%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  -0.04%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  +4.10%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '8.97'
$ws.Cells.Item(48, 5).Value = '  +10.33%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  +2.50%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '3.13'
$ws.Cells.Item(50, 5).Value = '  +1.38%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'Monero'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '143.42'
$ws.Cells.Item(51, 5).Value = '  +6.23%  '
